# Rename header labels on the first row of each worksheet so that when the
# table is imported into Power BI, the first row can be turned into headers
# automatically (labels must not be pure numbers like "2015").

$wb = $excel.ActiveWorkbook

# Sheets that use the "Ano <year>" pattern for B1:E1 (years 2015/2030/2040/2050)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Value2
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Value2
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Value2
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Value2
}

# Sheet that uses the "Intervalo <range>" pattern for B1:E1
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("B1").Value = "Intervalo " + $ws.Range("B1").Value2
$ws.Range("C1").Value = "Intervalo " + $ws.Range("C1").Value2
$ws.Range("D1").Value = "Intervalo " + $ws.Range("D1").Value2
$ws.Range("E1").Value = "Intervalo " + $ws.Range("E1").Value2

# Sheet with only a B1 header to update (no C1/D1/E1 columns present)
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Range("B1").Value = "Ano " + $ws.Range("B1").Value2
